# "regenerate orders with updates distance/sizes"
#
# The workbook encodes a trial-order table. Several distance codes and one
# size code were renumbered throughout the sheet (condition names, left/right
# filenames, and the standalone Distance/Size lookup columns all reference
# the same tokens):
#
#   D51 -> D55
#   D64 -> D69
#   D80 -> D86
#   S30 -> S31
#
# Every other token (S20, S25, Face##, etc.) is left untouched. Apply this as
# a whole-sheet, whole-cell-content text substitution so it touches every
# occurrence (Condition, Filename_Left, Filename_Right, Distance, Size
# columns) regardless of row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$replacements = @(
    @("D51", "D55"),
    @("D64", "D69"),
    @("D80", "D86"),
    @("S30", "S31")
)

foreach ($pair in $replacements) {
    $find = $pair[0]
    $replace = $pair[1]
    # xlWhole? No - tokens are embedded inside larger strings like
    # "Face03_D51_S25" or "Face03_D51_S25_l.png", so we need a substring
    # (xlPart = 2) replace, not a whole-cell match.
    $ws.Cells.Replace($find, $replace, 2)
}
